# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) timestamps on the per-language
# report sheets, reflecting a freshly regenerated handback report.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-19 06:50:51"
$zhcn.Range("H2").Value = "2016-03-19 06:51:16"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-19 06:50:54"
$dede.Range("H2").Value = "2016-03-19 06:51:22"
